$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 is missing a value in column C (faculrod) for the "Запасний" / "ЗАПАС" entry.
# Fill it in with the same text used in column D ("ЗАПАС"), matching the style
# already used by the sibling cells B15 and D15 in that row.
$ws.Range("C15").Value = "ЗАПАС"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Leave the selection where the editor ended up after filling the cell.
$ws.Range("F15").Select() | Out-Null
